# [external commands] - [tail(id,file)]: simulate the *NIX tail command.
#
# The "#system" sheet (sheet1) backs the data-validation dropdowns used on
# the visible sheets via named ranges: "external" (column I) lists the
# external-command signatures, "web" (column Y) lists the web-command
# signatures. Adding the new `tail(id,file)` external command appends a
# row to the "external" list; adding the new
# `assertTextNotContains(locator,text)` web command inserts a row into the
# alphabetically sorted "web" list (between assertTextMatches and
# assertTextNotPresent), pushing every subsequent entry down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. external commands (column I): append tail(id,file) at I5 ---------
$ws.Range("I5").Value2 = "tail(id,file)"

# --- 2. web commands (column Y): insert assertTextNotContains(locator,text)
#        at row 39, shifting Y39:Y127 down to Y40:Y128 -------------------
for ($r = 127; $r -ge 39; $r--) {
    $srcCell = $ws.Range("Y" + $r)
    $dstRow = $r + 1
    $dstCell = $ws.Range("Y" + $dstRow)
    $dstCell.Value2 = $srcCell.Value2
}
$ws.Range("Y39").Value2 = "assertTextNotContains(locator,text)"

# --- 3. keep the named ranges in sync with the new list extents ---------
$wb.Names.Item("external").RefersTo = "='#system'!`$I`$2:`$I`$5"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$128"
